$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "E2" = 3
    "F2" = 1
    "G2" = 4.202518666666666
    "H2" = 12.607556
    "I2" = 0.08075097102331126
    "J2" = 0.08075097102331129
    "M2" = 36.89194233333333
    "N2" = 110.675827
    "O2" = 0.3567095043190808
    "P2" = 0.3567095043190809
    "Q2" = 155.0390763054235
    "R2" = 1395.351686748812
    "S2" = 0.02880463884700982
    "T2" = 0.02880463884700983
    "E3" = 3
    "F3" = 1
    "G3" = 4.202518666666666
    "H3" = 12.607556
    "I3" = 0.08075097102331126
    "J3" = 0.08075097102331129
    "M3" = 42.68037399999999
    "O3" = 0.4126780562577495
    "P3" = 0.4126780562577496
    "Q3" = 179.3650684353146
    "R3" = 1614.285615917832
    "S3" = 0.03332415376282594
    "T3" = 0.03332415376282596
    "E4" = 3
    "F4" = 1
    "G4" = 4.202518666666666
    "H4" = 12.607556
    "I4" = 0.08075097102331126
    "J4" = 0.08075097102331129
    "M4" = 23.85061433333334
    "N4" = 71.55184300000001
    "O4" = 0.2306124394231696
    "P4" = 0.2306124394231696
    "Q4" = 100.2326519473009
    "R4" = 902.093867525708
    "S4" = 0.01862217841347549
    "T4" = 0.0186221784134755
    "I5" = 0.7496282157262072
    "J5" = 0.7496282157262073
    "M5" = 36.89194233333333
    "N5" = 110.675827
    "O5" = 0.3567095043190808
    "P5" = 0.3567095043190809
    "Q5" = 1439.260292054234
    "R5" = 12953.3426284881
    "S5" = 0.2673995092552924
    "T5" = 0.2673995092552924
    "I6" = 0.7496282157262072
    "J6" = 0.7496282157262073
    "M6" = 42.68037399999999
    "O6" = 0.4126780562577495
    "P6" = 0.4126780562577496
    "S6" = 0.3093551149818561
    "T6" = 0.3093551149818562
    "I7" = 0.7496282157262072
    "J7" = 0.7496282157262073
    "M7" = 23.85061433333334
    "N7" = 71.55184300000001
    "O7" = 0.2306124394231696
    "P7" = 0.2306124394231696
    "Q7" = 930.4807494521699
    "R7" = 8374.326745069528
    "S7" = 0.1728735914890586
    "T7" = 0.1728735914890587
    "G8" = 8.827567333333333
    "I8" = 0.1696208132504815
    "J8" = 0.1696208132504815
    "M8" = 36.89194233333333
    "N8" = 110.675827
    "O8" = 0.3567095043190808
    "P8" = 0.3567095043190809
    "Q8" = 325.6661050049504
    "R8" = 2930.994945044554
    "S8" = 0.06050535621677863
    "T8" = 0.06050535621677865
    "G9" = 8.827567333333333
    "I9" = 0.1696208132504815
    "J9" = 0.1696208132504815
    "M9" = 42.68037399999999
    "O9" = 0.4126780562577495
    "P9" = 0.4126780562577496
    "Q9" = 376.7638752968492
    "R9" = 3390.874877671643
    "S9" = 0.06999878751306741
    "T9" = 0.06999878751306746
    "G10" = 8.827567333333333
    "I10" = 0.1696208132504815
    "J10" = 0.1696208132504815
    "M10" = 23.85061433333334
    "N10" = 71.55184300000001
    "O10" = 0.2306124394231696
    "P10" = 0.2306124394231696
    "S10" = 0.03911666952063542
    "T10" = 0.03911666952063544
}

foreach ($cell in $values.Keys) {
    $ws.Range($cell).Value = $values[$cell]
}
